$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 from "fam_plan_need" to "fam_plan_unmet"
$ws.Range("E1").Value = "fam_plan_unmet"

# Update the selection to D17 (cosmetic, matches author's last selection before save)
$ws.Range("D17").Select()

# Adjust column E width slightly (matches autofit-style resize seen in diff)
$ws.Columns.Item(5).ColumnWidth = 14
